$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: "scale" row (5/10/15 repeated), styled like row 1's column groups ---
$ws.Range("B1:D1").Copy()
$ws.Range("B6:D6").PasteSpecial(-4122)
$ws.Range("E1:G1").Copy()
$ws.Range("E6:G6").PasteSpecial(-4122)
$ws.Range("H1:J1").Copy()
$ws.Range("H6:J6").PasteSpecial(-4122)
$ws.Range("K1:M1").Copy()
$ws.Range("K6:M6").PasteSpecial(-4122)
$ws.Range("N3:P3").Copy()
$ws.Range("N6:P6").PasteSpecial(-4122)

$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 15
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 15
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 10
$ws.Range("J6").Value = 15
$ws.Range("K6").Value = 5
$ws.Range("L6").Value = 10
$ws.Range("M6").Value = 15
$ws.Range("N6").Value = 5
$ws.Range("O6").Value = 10
$ws.Range("P6").Value = 15

# --- Row 5: blank spacer cell with centered alignment ---
$ws.Range("H5").HorizontalAlignment = -4108
$ws.Range("H5").VerticalAlignment = -4108

# --- Row 7: "ExternalSort1" results row, styled like rows 3/4 (right aligned fills) ---
$ws.Range("C3:D3").Copy()
$ws.Range("B7:D7").PasteSpecial(-4122)
$ws.Range("E3:G3").Copy()
$ws.Range("E7:G7").PasteSpecial(-4122)
$ws.Range("H3:J3").Copy()
$ws.Range("H7:J7").PasteSpecial(-4122)
$ws.Range("K4:M4").Copy()
$ws.Range("K7:M7").PasteSpecial(-4122)
$ws.Range("N4:P4").Copy()
$ws.Range("N7:P7").PasteSpecial(-4122)

$ws.Range("A7").Value = "ExternalSort1"
$ws.Range("B7").Value = "0.007"
$ws.Range("C7").Value = "0.009"
$ws.Range("D7").Value = "0.014"
$ws.Range("E7").Value = "0.011"
$ws.Range("F7").Value = "0.015"
$ws.Range("G7").Value = "0.020"
$ws.Range("H7").Value = "0.072"
$ws.Range("I7").Value = "0.86"
$ws.Range("J7").Value = "0.100"
$ws.Range("K7").Value = "0.766"
$ws.Range("L7").Value = "0.838"
$ws.Range("M7").Value = "0.924"
$ws.Range("N7").Value = "8.453"
$ws.Range("O7").Value = "9.232"
$ws.Range("P7").Value = "10.119"

# --- Column A width, view state ---
$ws.Columns("A").ColumnWidth = 12.77734375
$ws.Range("N9").Select()
